$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# New log entries for the mux8x1 work (rows 10-12), following on from the
# existing mux4x1 entries above them.
# ---------------------------------------------------------------------

# Row 10 - Sr No. 4 (B10 already = 4)
$ws.Range("C10").Value = 44879
$ws.Range("C10").NumberFormat = "d-mmm-yy"
$ws.Range("D10").Value = "mux8x1"
$ws.Range("E10").Value = "mux8x1.v"
$ws.Range("F10").Value = "mux8x1.v"
$ws.Range("G10").Value = "Warning: `"Port size (1) does not match connection size (32) for port i2`""
$ws.Range("H10").Value = "While assiging 0 to any variable use e.g. 1'b0 (one bit zero). By default, it will consider 32 bit 0 instead of your required bit zero and it will give you this error."

# Row 11 - continuation row (B11 stays blank)
$ws.Range("C11").Value = 44879
$ws.Range("C11").NumberFormat = "d-mmm-yy"
$ws.Range("D11").Value = "mux8x1"
$ws.Range("E11").Value = "mux8x1.v"
$ws.Range("F11").Value = "mux8x1.v"
$ws.Range("G11").Value = "Error (suppressible): (vsim-3053) Illegal output or inout port connection for port 'out'.`nTried to connect `"reg out (8x1) port`" to `"reg out (4x1) port`"."
$ws.Range("H11").Value = "port out (8x1) must be defined as wire as it is just driving the out (4x1) port."

# Row 12 - Sr No. 5 (B12 already = 5)
$ws.Range("C12").Value = 44879
$ws.Range("C12").NumberFormat = "d-mmm-yy"
$ws.Range("D12").Value = "mux8x1"
$ws.Range("E12").Value = "mux8x1_tb.v"
$ws.Range("F12").Value = "mux8x1_tb.v"
$ws.Range("G12").Value = "GetModuleFileName: The specified module could not be found"
$ws.Range("H12").Value = "Need to include mux4x1.v to the tb file.`nMux4x1 module is instantiated to mux8x1 and included in mux8x1.v file but not to the mux8x1_tb.v. So, by including 4x1 file to tb the error has been resolved."

# Let row heights fall back to the sheet's auto-computed height instead of
# leaving a stray "customHeight" flag behind from setting multi-line text.
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()

# ---------------------------------------------------------------------
# Normalise the Problem/Solution columns (G:H) so every row (including the
# new ones above) shares the same look: left aligned, vertically centred,
# wrapped text, regular (non-bold) font.
# ---------------------------------------------------------------------
$fmt = $ws.Range("G6:H39")
$fmt.HorizontalAlignment = -4131
$fmt.VerticalAlignment = -4108
$fmt.WrapText = $true
$fmt.Font.Bold = $false

# ---------------------------------------------------------------------
# Extend the same look down to row 92 for future entries, copying the
# exact formatting (border, alignment, wrap) from an already-normalised
# cell so the new blank rows match pixel for pixel.
# ---------------------------------------------------------------------
$ws.Range("G6:H6").Copy()
$ws.Range("G40:H92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Restore the active selection to where the author left off.
# ---------------------------------------------------------------------
$ws.Range("H13").Select()
